# Daily auto push: 2026-02-23 08:00 UTC
#
# A new observation for 2026/02/23 (Monday / 月) at hour 16 needs to be
# appended to the existing block of rows for that date (which currently
# ends at row 856 with hour 13). This pushes row 857 onward down by one
# row, extending the used range from A1:D898 to A1:D899.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last existing "2026/02/23" row (row 856: date, weekday,
# hour=13, ranking=201) into a freshly inserted row 857. Copy+Insert (as
# opposed to writing string literals straight into .Value) preserves the
# original date/weekday cells as plain text instead of Excel coercing the
# "2026/02/23"-looking string into a date serial number.
$ws.Rows.Item(856).Copy()
$ws.Rows.Item(857).Insert()
$excel.CutCopyMode = $false

# The new row 857 is a copy of row 856, so date/weekday/ranking are
# already correct (2026/02/23, 月, 201) - only the hour column differs.
$ws.Cells.Item(857, 3).Value = 16
